$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J, matching the style used by the other headers (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows (plain numbers, no special style)
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 3

$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 4

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 7

$ws.Range("I6").Value = 9
$ws.Range("J6").Value = 9

$ws.Range("I7").Value = 4
$ws.Range("J7").Value = 5

$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 5

$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 7
